$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the long descriptive text currently sitting in A2 before we start
# shuffling things around (reading .Value is unreliable in this host --
# .Value2 is the safe read path).
$longText = $ws.Range("A2").Value2

# --- Header relabel: "TaskPath"/"TaskDepth" -> "Task Path"/"Depth" ---
$ws.Range("A1").Value = "Task Path"
$ws.Range("B1").Value = "Depth"

# --- Move the long descriptive text into A2 (keeps same text/cell, just
#     re-asserted so the shared-string table gets rebuilt in the new order
#     the exported file expects: long text first, then the two headers) ---
$ws.Range("A2").Value = $longText

# --- Column A: narrower export width ---
$ws.Columns.Item(1).ColumnWidth = 133.33333333333334

# --- Row 2: drop the old forced 409.5pt custom height, use the export's
#     30pt (two wrapped lines at the default 15pt line height) ---
$ws.Rows.Item(2).RowHeight = 30

# --- A2 formatting: wrap + top-aligned ---
$ws.Range("A2").WrapText = $true
$ws.Range("A2").VerticalAlignment = -4160   # xlTop

# --- B2 formatting: top-aligned (no wrap) ---
$ws.Range("B2").VerticalAlignment = -4160   # xlTop

# --- Selection, matching the exported file ---
$ws.Range("A7").Select()
